# Apply the dial value corrections and view-state changes described by the diff.

$wb = $excel.ActiveWorkbook

$stacr  = $wb.Worksheets.Item("STACR")
$cas    = $wb.Worksheets.Item("CAS")
$nonqm  = $wb.Worksheets.Item("NONQM")

# --- STACR: correct the V2 Dial value for row 30 ---
$stacr.Range("H30").Value = 1.75

# --- CAS: correct the V2 Dial values for rows 22 and 28 ---
$cas.Range("H22").Value = 1.75
$cas.Range("H28").Value = 1.75

# --- NONQM: refresh V2 Dial values (tiny precision refresh + two real corrections) ---
$nonqm.Range("H29").Value = 1.09064438292682
$nonqm.Range("H30").Value = 1.1651372000584399
$nonqm.Range("H31").Value = 1.3456222354051699
$nonqm.Range("H32").Value = 1.60625232365774
$nonqm.Range("H33").Value = 0.73706281062687995
$nonqm.Range("H36").Value = 0.25
$nonqm.Range("H37").Value = 0.25
$nonqm.Range("H38").Value = 0.48833040506159497
$nonqm.Range("H41").Value = 1.0957770660791699
$nonqm.Range("H42").Value = 0.25

# --- Update sheet selections to match the saved cursor positions ---
[void]$stacr.Range("I23").Select()
[void]$cas.Range("A28").Select()
[void]$nonqm.Range("A44").Select()

# --- Make CAS the active sheet/tab (was STACR) ---
[void]$cas.Activate()
